$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update F-column "time_taken" timestamps on the data sheet ---
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:20:05.749985"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:20:05.749994"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:20:05.749997"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:20:05.750000"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:20:05.750003"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:20:05.750006"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:20:05.750009"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:20:05.750012"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:20:05.750014"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:20:05.750017"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:20:05.750020"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:20:05.750023"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:20:05.750025"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:20:05.750028"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:20:05.750031"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:20:05.750033"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:20:05.750036"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:20:05.750039"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:20:05.750042"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:20:05.750045"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:20:05.750047"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:20:05.750050"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:20:05.750053"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:20:05.750055"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:20:05.750059"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:20:05.750061"
$ws.Cells.Item(28, 6).Value = "2021-10-05 14:20:05.750064"
$ws.Cells.Item(29, 6).Value = "2021-10-05 14:20:05.750067"
$ws.Cells.Item(30, 6).Value = "2021-10-05 14:20:05.750069"
$ws.Cells.Item(31, 6).Value = "2021-10-05 14:20:05.750072"
$ws.Cells.Item(32, 6).Value = "2021-10-05 14:20:05.750075"
$ws.Cells.Item(33, 6).Value = "2021-10-05 14:20:05.750077"
$ws.Cells.Item(34, 6).Value = "2021-10-05 14:20:05.750081"
$ws.Cells.Item(35, 6).Value = "2021-10-05 14:20:05.750083"
$ws.Cells.Item(36, 6).Value = "2021-10-05 14:20:05.750086"
$ws.Cells.Item(37, 6).Value = "2021-10-05 14:20:05.750089"
$ws.Cells.Item(38, 6).Value = "2021-10-05 14:20:05.750092"
$ws.Cells.Item(39, 6).Value = "2021-10-05 14:20:05.750095"
$ws.Cells.Item(40, 6).Value = "2021-10-05 14:20:05.750097"
$ws.Cells.Item(41, 6).Value = "2021-10-05 14:20:05.750100"
$ws.Cells.Item(42, 6).Value = "2021-10-05 14:20:05.750103"
$ws.Cells.Item(43, 6).Value = "2021-10-05 14:20:05.750106"
$ws.Cells.Item(44, 6).Value = "2021-10-05 14:20:05.750109"
$ws.Cells.Item(45, 6).Value = "2021-10-05 14:20:05.750112"
$ws.Cells.Item(46, 6).Value = "2021-10-05 14:20:05.750115"
$ws.Cells.Item(47, 6).Value = "2021-10-05 14:20:05.750117"
$ws.Cells.Item(48, 6).Value = "2021-10-05 14:20:05.750120"
$ws.Cells.Item(49, 6).Value = "2021-10-05 14:20:05.750123"
$ws.Cells.Item(50, 6).Value = "2021-10-05 14:20:05.750126"
$ws.Cells.Item(51, 6).Value = "2021-10-05 14:20:05.750128"
$ws.Cells.Item(52, 6).Value = "2021-10-05 14:20:05.750131"
$ws.Cells.Item(53, 6).Value = "2021-10-05 14:20:05.750134"
$ws.Cells.Item(54, 6).Value = "2021-10-05 14:20:05.750137"
$ws.Cells.Item(55, 6).Value = "2021-10-05 14:20:05.750140"
$ws.Cells.Item(56, 6).Value = "2021-10-05 14:20:05.750143"
$ws.Cells.Item(57, 6).Value = "2021-10-05 14:20:05.750146"
$ws.Cells.Item(58, 6).Value = "2021-10-05 14:20:05.750148"
$ws.Cells.Item(59, 6).Value = "2021-10-05 14:20:05.750151"
$ws.Cells.Item(60, 6).Value = "2021-10-05 14:20:05.750154"
$ws.Cells.Item(61, 6).Value = "2021-10-05 14:20:05.750157"
$ws.Cells.Item(62, 6).Value = "2021-10-05 14:20:05.750160"
$ws.Cells.Item(63, 6).Value = "2021-10-05 14:20:05.750162"
$ws.Cells.Item(64, 6).Value = "2021-10-05 14:20:05.750165"
$ws.Cells.Item(65, 6).Value = "2021-10-05 14:20:05.750168"
$ws.Cells.Item(66, 6).Value = "2021-10-05 14:20:05.750172"
$ws.Cells.Item(67, 6).Value = "2021-10-05 14:20:05.750175"
$ws.Cells.Item(68, 6).Value = "2021-10-05 14:20:05.750179"
$ws.Cells.Item(69, 6).Value = "2021-10-05 14:20:05.750181"
$ws.Cells.Item(70, 6).Value = "2021-10-05 14:20:05.750184"
$ws.Cells.Item(71, 6).Value = "2021-10-05 14:20:05.750187"
$ws.Cells.Item(72, 6).Value = "2021-10-05 14:20:05.750190"
$ws.Cells.Item(73, 6).Value = "2021-10-05 14:20:05.750193"
$ws.Cells.Item(74, 6).Value = "2021-10-05 14:20:05.750196"
$ws.Cells.Item(75, 6).Value = "2021-10-05 14:20:05.750199"
$ws.Cells.Item(76, 6).Value = "2021-10-05 14:20:05.750202"
$ws.Cells.Item(77, 6).Value = "2021-10-05 14:20:05.750205"

# --- Add the new "metadata" sheet directly after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Ectodermal dysplasia"
$meta.Cells.Item(2, 3).Value = 553
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.26"
$meta.Cells.Item(2, 4).Style = "Normal"
$meta.Cells.Item(2, 5).Value = "2021-07-28T14:12:26.226916Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:05.746421"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/553/?format=json"

# Mirror the "data" sheet's header/index formatting (bold, centered, bordered)
# onto the new sheet's header row and A2 index cell, without introducing
# any new style entries.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep "data" as the active/selected sheet (matches the original workbook's
# activeTab, which the diff does not change).
$ws.Activate()
